# The ghosts' jail (columns L-Q, rows 14-16) was previously a block of
# walls (value 0). Open it up into empty/pellet tiles (value 2) so the
# ghosts are able to move and break out when returning home.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L14:Q14").Value = 2
$ws.Range("L15:Q15").Value = 2
$ws.Range("L16:Q16").Value = 2

# Leave the active selection on L16 (previously it was on Q16).
$ws.Range("L16").Select()
